# Updated cryptos list on Sat Jan 27 21:25:39 UTC 2024 with GitHub Actions
# Refresh the Price (D) / Volume(1h) (E) figures, and re-rank a block of
# coins (B/C/D/E) on rows 24-48 to reflect the latest coinranking.com pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference an already-unstyled data cell so that re-applying text values
# (below) never leaves a stray "quote prefix" / number-format style behind
# on the edited cells.
$refStyle = $ws.Range('B2').Style

$updates = @(
    @{Cell='D2'; Value='''42.043.34'},
    @{Cell='D3'; Value='''2.269.93'},
    @{Cell='E3'; Value='''  +0.69%  '},
    @{Cell='E4'; Value='''  -0.01%  '},
    @{Cell='D5'; Value='''305.73'},
    @{Cell='E5'; Value='''  +1.36%  '},
    @{Cell='D6'; Value='''93.24'},
    @{Cell='E6'; Value='''  +1.25%  '},
    @{Cell='E7'; Value='''  -0.20%  '},
    @{Cell='E8'; Value='''  -0.06%  '},
    @{Cell='E9'; Value='''  +1.50%  '},
    @{Cell='D10'; Value='''32.82'},
    @{Cell='E11'; Value='''  +0.27%  '},
    @{Cell='E12'; Value='''  -1.72%  '},
    @{Cell='D13'; Value='''6.69'},
    @{Cell='E13'; Value='''  +0.52%  '},
    @{Cell='D14'; Value='''2.620.44'},
    @{Cell='E14'; Value='''  +0.65%  '},
    @{Cell='D15'; Value='''14.35'},
    @{Cell='E15'; Value='''  +1.90%  '},
    @{Cell='D16'; Value='''2.271.35'},
    @{Cell='E16'; Value='''  +0.43%  '},
    @{Cell='D17'; Value='''0.786'},
    @{Cell='E17'; Value='''  +3.91%  '},
    @{Cell='D18'; Value='''41.918.67'},
    @{Cell='E18'; Value='''  +0.33%  '},
    @{Cell='D19'; Value='''12.78'},
    @{Cell='E19'; Value='''  +5.52%  '},
    @{Cell='E20'; Value='''  +1.84%  '},
    @{Cell='E21'; Value='''  +1.37%  '},
    @{Cell='D22'; Value='''68.17'},
    @{Cell='E22'; Value='''  +1.80%  '},
    @{Cell='D23'; Value='''244.38'},
    @{Cell='E23'; Value='''  +1.42%  '},
    @{Cell='B24'; Value='BinanceUSD'},
    @{Cell='C24'; Value='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'},
    @{Cell='D24'; Value='''8.41'},
    @{Cell='E24'; Value='''  +740.58%  '},
    @{Cell='B25'; Value='PancakeSwap'},
    @{Cell='C25'; Value='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'},
    @{Cell='D25'; Value='''2.62'},
    @{Cell='E25'; Value='''  +2.46%  '},
    @{Cell='B26'; Value='ImmutableX'},
    @{Cell='C26'; Value='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'},
    @{Cell='D26'; Value='''1.94'},
    @{Cell='E26'; Value='''  +2.61%  '},
    @{Cell='B27'; Value='Dai'},
    @{Cell='C27'; Value='https://coinranking.com/coin/MoTuySvg7+dai-dai'},
    @{Cell='D27'; Value='''1.00'},
    @{Cell='E27'; Value='''  +0.06%  '},
    @{Cell='B28'; Value='EthereumClassic'},
    @{Cell='C28'; Value='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'},
    @{Cell='D28'; Value='''24.02'},
    @{Cell='E28'; Value='''  +0.52%  '},
    @{Cell='B29'; Value='Cosmos'},
    @{Cell='C29'; Value='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'},
    @{Cell='D29'; Value='''9.68'},
    @{Cell='E29'; Value='''  +0.40%  '},
    @{Cell='B30'; Value='Toncoin'},
    @{Cell='C30'; Value='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'},
    @{Cell='D30'; Value='''2.09'},
    @{Cell='E30'; Value='''  -9.34%  '},
    @{Cell='B31'; Value='InjectiveProtocol'},
    @{Cell='C31'; Value='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'},
    @{Cell='D31'; Value='''35.04'},
    @{Cell='E31'; Value='''  +3.90%  '},
    @{Cell='B32'; Value='Monero'},
    @{Cell='C32'; Value='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'},
    @{Cell='D32'; Value='''159.77'},
    @{Cell='E32'; Value='''  +0.63%  '},
    @{Cell='B33'; Value='Filecoin'},
    @{Cell='C33'; Value='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'},
    @{Cell='D33'; Value='''5.35'},
    @{Cell='E33'; Value='''  +3.99%  '},
    @{Cell='B34'; Value='FirstDigitalUSD'},
    @{Cell='C34'; Value='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'},
    @{Cell='D34'; Value='''0.999'},
    @{Cell='E34'; Value='''  +0.02%  '},
    @{Cell='B35'; Value='Hedera'},
    @{Cell='C35'; Value='https://coinranking.com/coin/jad286TjB+hedera-hbar'},
    @{Cell='D35'; Value='''0.0744'},
    @{Cell='E35'; Value='''  +0.30%  '},
    @{Cell='B36'; Value='LidoDAOToken'},
    @{Cell='C36'; Value='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'},
    @{Cell='D36'; Value='''3.03'},
    @{Cell='E36'; Value='''  -0.14%  '},
    @{Cell='B37'; Value='Celestia'},
    @{Cell='C37'; Value='https://coinranking.com/coin/YQcD0lBl7+celestia-tia'},
    @{Cell='D37'; Value='''17.17'},
    @{Cell='E37'; Value='''  +4.60%  '},
    @{Cell='B38'; Value='WEMIXToken'},
    @{Cell='C38'; Value='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'},
    @{Cell='D38'; Value='''2.37'},
    @{Cell='E38'; Value='''  -1.08%  '},
    @{Cell='B39'; Value='Kaspa'},
    @{Cell='C39'; Value='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'},
    @{Cell='D39'; Value='''0.105'},
    @{Cell='E39'; Value='''  +1.45%  '},
    @{Cell='B40'; Value='Stellar'},
    @{Cell='C40'; Value='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'},
    @{Cell='D40'; Value='''0.117'},
    @{Cell='E40'; Value='''  +1.07%  '},
    @{Cell='B41'; Value='ARBITRUM'},
    @{Cell='C41'; Value='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'},
    @{Cell='D41'; Value='''1.80'},
    @{Cell='E41'; Value='''  +0.67%  '},
    @{Cell='B42'; Value='RenderToken'},
    @{Cell='C42'; Value='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'},
    @{Cell='D42'; Value='''4.00'},
    @{Cell='E42'; Value='''  +2.00%  '},
    @{Cell='B43'; Value='EnergySwap'},
    @{Cell='C43'; Value='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'},
    @{Cell='D43'; Value='''19.85'},
    @{Cell='E43'; Value='''  +1.02%  '},
    @{Cell='B44'; Value='Maker'},
    @{Cell='C44'; Value='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'},
    @{Cell='D44'; Value='''2.017.06'},
    @{Cell='E44'; Value='''  -1.53%  '},
    @{Cell='B45'; Value='ApeXProtocol'},
    @{Cell='C45'; Value='https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'},
    @{Cell='D45'; Value='''2.24'},
    @{Cell='E45'; Value='''  +9.41%  '},
    @{Cell='B46'; Value='VeChain'},
    @{Cell='C46'; Value='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'},
    @{Cell='D46'; Value='''0.0283'},
    @{Cell='E46'; Value='''  +1.50%  '},
    @{Cell='B47'; Value='FraxShare'},
    @{Cell='C47'; Value='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'},
    @{Cell='D47'; Value='''10.26'},
    @{Cell='E47'; Value='''  +2.05%  '},
    @{Cell='B48'; Value='NEARProtocol'},
    @{Cell='C48'; Value='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'},
    @{Cell='D48'; Value='''2.91'},
    @{Cell='E48'; Value='''  +2.29%  '},
    @{Cell='D49'; Value='''53.34'},
    @{Cell='E49'; Value='''  +3.48%  '},
    @{Cell='E50'; Value='''  +0.53%  '},
    @{Cell='D51'; Value='''72.50'},
    @{Cell='E51'; Value='''  +2.91%  '}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.Value2 = $u.Value
    $cell.Style = $refStyle
}
